$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''63.809.61'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '''3.322.76'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''603.71'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').Value = '''142.93'
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '''3.320.18'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('D11').Value = '''5.55'
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').Value = '''0.0000248'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').Value = '''35.04'
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = '''3.868.38'
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('D16').Value = '''0.121'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '''3.323.33'
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('D18').Value = '''63.874.16'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('E19').Value = '  +1.31%  '
$ws.Range('D20').Value = '''480.67'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('E22').Value = '  +2.04%  '
$ws.Range('D23').Value = '''7.97'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('E24').Value = '  +4.05%  '
$ws.Range('D25').Value = '''84.84'
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +1.84%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '''8.27'
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').Value = '''7.21'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  +1.30%  '
$ws.Range('D32').Value = '''28.96'
$ws.Range('E32').Value = '  +4.99%  '
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('E35').Value = '  +1.29%  '
$ws.Range('E36').Value = '  +3.26%  '
$ws.Range('D37').Value = '''0.0₃0749'
$ws.Range('E37').Value = '  +5.03%  '
$ws.Range('D38').Value = '''52.39'
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('D39').Value = '''0.0399'
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '''3.116.25'
$ws.Range('E40').Value = '  +4.18%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '''431.54'
$ws.Range('E41').Value = '  +2.71%  '
$ws.Range('D42').Value = '''0.117'
$ws.Range('E42').Value = '  +5.32%  '
$ws.Range('B43').Value = 'Cosmos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D43').Value = '''8.35'
$ws.Range('E43').Value = '  -0.49%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '''2.75'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('D46').Value = '''2.24'
$ws.Range('E46').Value = '  +3.76%  '
$ws.Range('D47').Value = '''36.46'
$ws.Range('E47').Value = '  +9.38%  '
$ws.Range('D48').Value = '''26.42'
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.114'
$ws.Range('E51').Value = '  -0.56%  '
